$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45939.0
$ws.Cells.Item(2, 2).Value = 4282.70362407058
$ws.Cells.Item(2, 3).Value = 5028.80242723193
$ws.Cells.Item(2, 4).Value = 6812.0
$ws.Cells.Item(2, 5).Value = 5920.550514
$ws.Cells.Item(2, 6).Value = -6.05627845161041

$ws.Cells.Item(3, 1).Value = 45940.0
$ws.Cells.Item(3, 2).Value = 5212.44261081921
$ws.Cells.Item(3, 3).Value = 4597.10364210284
$ws.Cells.Item(3, 5).Value = 7525.272297
$ws.Cells.Item(3, 6).Value = 137.080555345151

$ws.Cells.Item(4, 1).Value = 45941.0
$ws.Cells.Item(4, 2).Value = 1359.17670681482
$ws.Cells.Item(4, 3).Value = 3057.29296888109
$ws.Cells.Item(4, 5).Value = 3241.094212
$ws.Cells.Item(4, 6).Value = 54.9671030860945

$ws.Cells.Item(5, 1).Value = 45942.0
$ws.Cells.Item(5, 2).Value = 1222.03028879225
$ws.Cells.Item(5, 3).Value = 3066.42463550278
$ws.Cells.Item(5, 5).Value = 3048.927259
$ws.Cells.Item(5, 6).Value = 53.0550669046055

$ws.Cells.Item(6, 1).Value = 45943.0
$ws.Cells.Item(6, 2).Value = 5417.79389491431
$ws.Cells.Item(6, 3).Value = 5318.07004356478
$ws.Cells.Item(6, 5).Value = 7930.914071
$ws.Cells.Item(6, 6).Value = 175.466259152103

$ws.Cells.Item(7, 1).Value = 45944.0
$ws.Cells.Item(7, 3).Value = 6012.25078158027
$ws.Cells.Item(7, 5).Value = 8951.45725
$ws.Cells.Item(7, 6).Value = 229.267678695247

$ws.Cells.Item(8, 1).Value = 45945.0
$ws.Cells.Item(8, 3).Value = 6406.23506180018
$ws.Cells.Item(8, 5).Value = 8971.008172
$ws.Cells.Item(8, 6).Value = 246.498312121076

$ws.Cells.Item(9, 1).Value = 45946.0
$ws.Cells.Item(9, 3).Value = 6185.47077103962
$ws.Cells.Item(9, 5).Value = 8971.008172
$ws.Cells.Item(9, 6).Value = 237.299800006053

$ws.Cells.Item(10, 1).Value = 45947.0
$ws.Cells.Item(10, 3).Value = 5031.6600158503
$ws.Cells.Item(10, 5).Value = 8971.008172
$ws.Cells.Item(10, 6).Value = 189.224351873165

$ws.Cells.Item(11, 1).Value = 45948.0
$ws.Cells.Item(11, 2).Value = 1742.27770790123
$ws.Cells.Item(11, 3).Value = 3457.00633017189
$ws.Cells.Item(11, 5).Value = 4473.853177
$ws.Cells.Item(11, 6).Value = 107.024241636277

$ws.Cells.Item(12, 1).Value = 45949.0
$ws.Cells.Item(12, 2).Value = 1636.94065696827
$ws.Cells.Item(12, 3).Value = 3436.36753946972
$ws.Cells.Item(12, 5).Value = 4358.692076
$ws.Cells.Item(12, 6).Value = 105.754956604227

$ws.Cells.Item(13, 1).Value = 45950.0
$ws.Cells.Item(13, 2).Value = 6392.95297294923
$ws.Cells.Item(13, 3).Value = 5964.83267309547
$ws.Cells.Item(13, 5).Value = 9928.522499
$ws.Cells.Item(13, 6).Value = 245.01675829776

$ws.Cells.Item(14, 1).Value = 45951.0
$ws.Cells.Item(14, 3).Value = 6197.27886190519
$ws.Cells.Item(14, 5).Value = 9928.522499
$ws.Cells.Item(14, 6).Value = 254.702016164832

$ws.Cells.Item(15, 1).Value = 45952.0
$ws.Cells.Item(15, 3).Value = 6437.54879871881
$ws.Cells.Item(15, 5).Value = 9928.522499
$ws.Cells.Item(15, 6).Value = 264.713263532066
